$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.989395499229431
$ws.Range("B1").Value = 2.270251750946045
$ws.Range("C1").Value = 2.268316507339478
$ws.Range("D1").Value = 2.734574794769287
$ws.Range("E1").Value = 3.507378101348877
